$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{ Cell = 'D2'; Value = '51.940.85' }
    @{ Cell = 'E2'; Value = '  +0.05%  ' }
    @{ Cell = 'D3'; Value = '2.789.39' }
    @{ Cell = 'E3'; Value = '  -1.67%  ' }
    @{ Cell = 'E4'; Value = '  -0.06%  ' }
    @{ Cell = 'D5'; Value = '358.62' }
    @{ Cell = 'E5'; Value = '  -0.04%  ' }
    @{ Cell = 'D6'; Value = '109.76' }
    @{ Cell = 'E6'; Value = '  -2.91%  ' }
    @{ Cell = 'D7'; Value = '0.560' }
    @{ Cell = 'E8'; Value = '  -0.03%  ' }
    @{ Cell = 'D9'; Value = '0.590' }
    @{ Cell = 'E9'; Value = '  -1.89%  ' }
    @{ Cell = 'D10'; Value = '40.43' }
    @{ Cell = 'E10'; Value = '  -2.26%  ' }
    @{ Cell = 'E11'; Value = '  +1.87%  ' }
    @{ Cell = 'D12'; Value = '0.0850' }
    @{ Cell = 'E12'; Value = '  -1.27%  ' }
    @{ Cell = 'D13'; Value = '19.51' }
    @{ Cell = 'E13'; Value = '  -3.23%  ' }
    @{ Cell = 'E14'; Value = '  -2.94%  ' }
    @{ Cell = 'D15'; Value = '3.228.29' }
    @{ Cell = 'E15'; Value = '  -1.93%  ' }
    @{ Cell = 'D16'; Value = '2.787.92' }
    @{ Cell = 'E16'; Value = '  -1.59%  ' }
    @{ Cell = 'D17'; Value = '0.954' }
    @{ Cell = 'E17'; Value = '  +2.79%  ' }
    @{ Cell = 'D18'; Value = '51.840.69' }
    @{ Cell = 'E18'; Value = '  -0.03%  ' }
    @{ Cell = 'D19'; Value = '7.46' }
    @{ Cell = 'E19'; Value = '  -1.46%  ' }
    @{ Cell = 'E20'; Value = '  -2.50%  ' }
    @{ Cell = 'D21'; Value = '13.19' }
    @{ Cell = 'E21'; Value = '  -2.24%  ' }
    @{ Cell = 'D22'; Value = '0.0₃0978' }
    @{ Cell = 'E22'; Value = '  -1.58%  ' }
    @{ Cell = 'D23'; Value = '270.99' }
    @{ Cell = 'E23'; Value = '  +0.55%  ' }
    @{ Cell = 'D24'; Value = '70.14' }
    @{ Cell = 'E24'; Value = '  -0.11%  ' }
    @{ Cell = 'E25'; Value = '  -2.72%  ' }
    @{ Cell = 'D26'; Value = '26.47' }
    @{ Cell = 'E26'; Value = '  -2.28%  ' }
    @{ Cell = 'E27'; Value = '  -0.02%  ' }
    @{ Cell = 'D28'; Value = '0.165' }
    @{ Cell = 'E28'; Value = '  +18.18%  ' }
    @{ Cell = 'D29'; Value = '10.28' }
    @{ Cell = 'E29'; Value = '  -0.73%  ' }
    @{ Cell = 'D30'; Value = '2.15' }
    @{ Cell = 'E30'; Value = '  -4.73%  ' }
    @{ Cell = 'B31'; Value = 'OKB' }
    @{ Cell = 'C31'; Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb' }
    @{ Cell = 'D31'; Value = '52.04' }
    @{ Cell = 'E31'; Value = '  -2.29%  ' }
    @{ Cell = 'B32'; Value = 'InjectiveProtocol' }
    @{ Cell = 'C32'; Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj' }
    @{ Cell = 'D32'; Value = '34.85' }
    @{ Cell = 'E32'; Value = '  -1.76%  ' }
    @{ Cell = 'B33'; Value = 'VeChain' }
    @{ Cell = 'C33'; Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet' }
    @{ Cell = 'D33'; Value = '0.0466' }
    @{ Cell = 'E33'; Value = '  -1.42%  ' }
    @{ Cell = 'D34'; Value = '5.76' }
    @{ Cell = 'E34'; Value = '  -2.77%  ' }
    @{ Cell = 'D35'; Value = '0.0848' }
    @{ Cell = 'E35'; Value = '  -1.14%  ' }
    @{ Cell = 'D36'; Value = '5.20' }
    @{ Cell = 'E36'; Value = '  -4.85%  ' }
    @{ Cell = 'E37'; Value = '  -0.09%  ' }
    @{ Cell = 'D38'; Value = '18.78' }
    @{ Cell = 'E38'; Value = '  +0.73%  ' }
    @{ Cell = 'D39'; Value = '3.21' }
    @{ Cell = 'E39'; Value = '  -2.69%  ' }
    @{ Cell = 'E40'; Value = '  -4.00%  ' }
    @{ Cell = 'D41'; Value = '2.58' }
    @{ Cell = 'E41'; Value = '  +1.35%  ' }
    @{ Cell = 'E42'; Value = '  -1.74%  ' }
    @{ Cell = 'E43'; Value = '  -1.94%  ' }
    @{ Cell = 'D44'; Value = '119.42' }
    @{ Cell = 'E44'; Value = '  -4.50%  ' }
    @{ Cell = 'D45'; Value = '21.81' }
    @{ Cell = 'E45'; Value = '  -7.43%  ' }
    @{ Cell = 'D46'; Value = '2.080.61' }
    @{ Cell = 'E46'; Value = '  -1.23%  ' }
    @{ Cell = 'D47'; Value = '3.28' }
    @{ Cell = 'E47'; Value = '  -3.33%  ' }
    @{ Cell = 'D48'; Value = '2.24' }
    @{ Cell = 'E48'; Value = '  -0.87%  ' }
    @{ Cell = 'D49'; Value = '5.79' }
    @{ Cell = 'E49'; Value = '  -3.46%  ' }
    @{ Cell = 'E50'; Value = '  -3.39%  ' }
    @{ Cell = 'B51'; Value = 'BitgetToken' }
    @{ Cell = 'C51'; Value = 'https://coinranking.com/coin/q7gMmMdLb+bitgettoken-bgb' }
    @{ Cell = 'D51'; Value = '1.13' }
    @{ Cell = 'E51'; Value = '  +29.81%  ' }
)

foreach ($chg in $changes) {
    $rng = $ws.Range($chg.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $chg.Value
    $rng.Style = "Normal"
}
